# regen sval data to filter save games
# Update the numeric stat columns (B..E, G) for rows 2-7 on the active sheet.
# Column F (Win) is left untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.445647641019636
    "C2" = 1.626987699542094
    "D2" = 3.223369029078222
    "E2" = 13.86384647080068
    "G2" = 20.15985084044064

    "B3" = 3.272327238179451
    "C3" = 1.626987699542094
    "D3" = 189.6080260415259
    "E3" = 0.5333859586016987
    "G3" = 195.0407269378492

    "B4" = 0.6545652718822623
    "C4" = 0.3048912486333797
    "D4" = 0.7210945179870265
    "E4" = 0.5333859586016987
    "G4" = 2.213936997104367

    "B5" = 3.272327238179451
    "C5" = 1.626987699542094
    "D5" = 0.7210945179870265
    "E5" = 0.5333859586016987
    "G5" = 6.15379541431027

    "B6" = 3.272327238179451
    "C6" = 1.626987699542094
    "D6" = 0.1496068669990043
    "E6" = 13.86384647080068
    "G6" = 18.91276827552123

    "B7" = 3.272327238179451
    "C7" = 1.626987699542094
    "D7" = 0.1496068669990043
    "E7" = 0.5333859586016987
    "G7" = 5.582307763322248
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
